# Update "想去人数" (interest count) figures on the 展览 and 全部类型 sheets,
# mirroring the latest scrape output.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 85
    $ws.Range("F4").Value = 1486
    $ws.Range("F9").Value = 271
}
